# Append new Lancers listing at row 6 (pushing existing rows 6-16 down to 7-17),
# refresh the "fetched at" timestamp for every data row, widen column H by one
# character, and rebuild the URL hyperlinks so they stay aligned with their rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-17 01:44:32"

# --- 1. Shift rows 6..16 down to 7..17 (bottom-up so we never clobber a row
#        before it has been copied). Column H is not included in the source
#        range for row 16 since it has no H value there; ClearContents below
#        removes any stray empty cell the copy leaves behind in that column
#        for rows that should not have one.
for ($r = 16; $r -ge 6; $r--) {
    $dstRow = $r + 1
    $srcRange = $ws.Range("A" + $r + ":H" + $r)
    $dstRange = $ws.Range("A" + $dstRow + ":H" + $dstRow)
    $srcRange.Copy($dstRange)
}

# Rows 16 and 17 (originally rows 15 and 16) never had a "skill summary" (H)
# value - drop the blank cell the block-copy above may have introduced there.
$ws.Range("H16").ClearContents()
$ws.Range("H17").ClearContents()

# --- 2. Write the brand-new row 6 ---
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5398112"
$ws.Range("G6").Value = 158
$ws.Range("H6").Value = "◆ツール,開発 ◇管理"

# --- 3. Refresh the timestamp column for every other data row (2-5 keep
#        their row position; 7-17 now hold what used to be rows 6-16) ---
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp
$ws.Range("A13").Value = $newTimestamp
$ws.Range("A14").Value = $newTimestamp
$ws.Range("A15").Value = $newTimestamp
$ws.Range("A16").Value = $newTimestamp
$ws.Range("A17").Value = $newTimestamp

# --- 4. Column H gets one character wider (12 -> 13 in raw OOXML units).
#        Excel's ColumnWidth property is offset from the stored <col width>
#        by ~5/6 of a character, so subtract that to land exactly on 13.
$ws.Columns.Item(8).ColumnWidth = 13 - (5/6)

# --- 5. Rebuild the hyperlinks for column F so ref cells line up with the
#        shifted rows. Per-item Hyperlink.Delete() is a no-op in this engine
#        and Range(...).Hyperlinks.Delete() clears the whole sheet regardless
#        of the range, so just clear everything once and re-add in order.
$ws.Hyperlinks.Delete()

# Keyed by the FINAL row number (post-insert) the URL belongs in.
$urls = @{
    2  = "https://www.lancers.jp/work/detail/5413954"
    3  = "https://www.lancers.jp/work/detail/5413955"
    4  = "https://www.lancers.jp/work/detail/5217096"
    5  = "https://www.lancers.jp/work/detail/5405023"
    6  = "https://www.lancers.jp/work/detail/5398112"
    7  = "https://www.lancers.jp/work/detail/5414167"
    8  = "https://www.lancers.jp/work/detail/5414354"
    9  = "https://www.lancers.jp/work/detail/5414353"
    10 = "https://www.lancers.jp/work/detail/5414105"
    11 = "https://www.lancers.jp/work/detail/5414108"
    12 = "https://www.lancers.jp/work/detail/5413916"
    13 = "https://www.lancers.jp/work/detail/5414368"
    14 = "https://www.lancers.jp/work/detail/5413958"
    15 = "https://www.lancers.jp/work/detail/5414569"
    16 = "https://www.lancers.jp/work/detail/5414579"
    17 = "https://www.lancers.jp/work/detail/5414812"
}

for ($rowNum = 2; $rowNum -le 17; $rowNum++) {
    $target = $urls[$rowNum]
    $cell = $ws.Range("F" + $rowNum)
    $cell.Value = $target
    $ws.Hyperlinks.Add($cell, $target)
    $cell.Style = "Hyperlink"
}
